# Update TPM-derived NATMI LR-pair metrics for Wnt16-Lrp5 (commit: "update scripts wuth new tpm")
# Only the "ECs" sending/target-cluster TPM-derived base values changed; this cascades
# through the specificity and edge-weight columns for every row that involves ECs
# (as sender and/or target), and produces tiny floating point re-normalisation noise
# in rows that only involve the unchanged clusters (FAPs/MuSCs) as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> ECs) ---
$ws.Range("G2").Value = 0.06577466666666666
$ws.Range("H2").Value = 0.197324
$ws.Range("I2").Value = 0.1530524945763335
$ws.Range("J2").Value = 0.1530524945763335
$ws.Range("M2").Value = 15.01856033333333
$ws.Range("N2").Value = 45.055681
$ws.Range("O2").Value = 0.4908713633047416
$ws.Range("P2").Value = 0.4908713633047417
$ws.Range("Q2").Value = 0.9878407997382221
$ws.Range("R2").Value = 8.890567197644
$ws.Range("S2").Value = 0.07512908666987639
$ws.Range("T2").Value = 0.0751290866698764

# --- Row 3 (ECs -> FAPs) ---
$ws.Range("G3").Value = 0.06577466666666666
$ws.Range("H3").Value = 0.197324
$ws.Range("I3").Value = 0.1530524945763335
$ws.Range("J3").Value = 0.1530524945763335
$ws.Range("O3").Value = 0.3099803572711625
$ws.Range("P3").Value = 0.3099803572711625
$ws.Range("Q3").Value = 0.623811586743111
$ws.Range("R3").Value = 5.614304280688
$ws.Range("S3").Value = 0.04744326695001451
$ws.Range("T3").Value = 0.04744326695001452

# --- Row 4 (ECs -> MuSCs) ---
$ws.Range("G4").Value = 0.06577466666666666
$ws.Range("H4").Value = 0.197324
$ws.Range("I4").Value = 0.1530524945763335
$ws.Range("J4").Value = 0.1530524945763335
$ws.Range("O4").Value = 0.1991482794240958
$ws.Range("P4").Value = 0.1991482794240958
$ws.Range("Q4").Value = 0.400770569072
$ws.Range("R4").Value = 3.606935121648
$ws.Range("S4").Value = 0.03048014095644257
$ws.Range("T4").Value = 0.03048014095644257

# --- Row 5 (FAPs -> ECs) ---
$ws.Range("H5").Value = 0.8943449999999999
$ws.Range("I5").Value = 0.6936902417438882
$ws.Range("J5").Value = 0.693690241743888
$ws.Range("M5").Value = 15.01856033333333
$ws.Range("N5").Value = 45.055681
$ws.Range("O5").Value = 0.4908713633047416
$ws.Range("P5").Value = 0.4908713633047417
$ws.Range("Q5").Value = 4.477258113771666
$ws.Range("R5").Value = 40.295323023945
$ws.Range("S5").Value = 0.3405126746760181
$ws.Range("T5").Value = 0.3405126746760181

# --- Row 6 (FAPs -> FAPs) ---
$ws.Range("H6").Value = 0.8943449999999999
$ws.Range("I6").Value = 0.6936902417438882
$ws.Range("J6").Value = 0.693690241743888
$ws.Range("O6").Value = 0.3099803572711625
$ws.Range("P6").Value = 0.3099803572711625
$ws.Range("S6").Value = 0.2150303489712895
$ws.Range("T6").Value = 0.2150303489712895

# --- Row 7 (FAPs -> MuSCs) ---
$ws.Range("H7").Value = 0.8943449999999999
$ws.Range("I7").Value = 0.6936902417438882
$ws.Range("J7").Value = 0.693690241743888
$ws.Range("O7").Value = 0.1991482794240958
$ws.Range("P7").Value = 0.1991482794240958
$ws.Range("S7").Value = 0.1381472180965804
$ws.Range("T7").Value = 0.1381472180965804

# --- Row 8 (MuSCs -> ECs) ---
$ws.Range("I8").Value = 0.1532572636797784
$ws.Range("J8").Value = 0.1532572636797783
$ws.Range("M8").Value = 15.01856033333333
$ws.Range("N8").Value = 45.055681
$ws.Range("O8").Value = 0.4908713633047416
$ws.Range("P8").Value = 0.4908713633047417
$ws.Range("Q8").Value = 0.9891624330475556
$ws.Range("R8").Value = 8.902461897428001
$ws.Range("S8").Value = 0.07522960195884708
$ws.Range("T8").Value = 0.07522960195884706

# --- Row 9 (MuSCs -> FAPs) ---
$ws.Range("I9").Value = 0.1532572636797784
$ws.Range("J9").Value = 0.1532572636797783
$ws.Range("O9").Value = 0.3099803572711625
$ws.Range("P9").Value = 0.3099803572711625
$ws.Range("S9").Value = 0.04750674134985845
$ws.Range("T9").Value = 0.04750674134985845

# --- Row 10 (MuSCs -> MuSCs) ---
$ws.Range("I10").Value = 0.1532572636797784
$ws.Range("J10").Value = 0.1532572636797783
$ws.Range("O10").Value = 0.1991482794240958
$ws.Range("P10").Value = 0.1991482794240958
$ws.Range("S10").Value = 0.03052092037107284
$ws.Range("T10").Value = 0.03052092037107283
